# "Wireframes version 2." -> "Wireframes version 1."
# The paragraph text is split across several runs:
#   "Versi" | "on" | " 2" | "."   (with a proofErr spellcheck wrapper
#   around "Versi"+"on" and a _GoBack bookmark after " 2")
# Target state merges "Versi"+"on" into a single "Version" run, changes
# " 2" to " 1.", and drops the now-redundant trailing "." run.

$d = $word.ActiveDocument

# 1) Merge the split "Versi" + "on" runs back into a single "Version" run.
$d.Content.Find.Execute("Version", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Version", 2)

# 2) Change the version number " 2" -> " 1." (the trailing period moves
#    here so the old standalone "." run becomes redundant).
$d.Content.Find.Execute(" 2", $true, $false, $false, $false, $false,
                         $true, 1, $false, " 1.", 2)

# 3) Remove the now-duplicate trailing "." run (the one after the
#    _GoBack bookmark), leaving "Version 1." as the full paragraph text.
$count = $d.Characters.Count
$last = $d.Range($count - 2, $count - 1)
if ($last.Text -eq ".") {
    $last.Delete()
}
